# Insert a new row at row 46, pushing existing rows 46-74 down to 47-75.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("46:46").Insert()

# Populate the newly inserted row 46 with the same data as the (now shifted)
# row 47 -- i.e. a duplicate of the original row 46 -- except for the new
# Fecha (D) and Volumen (J) values.
$ws.Range("A46").Value = 10
$ws.Range("B46").Value = "Vega Modelo de Temuco"
$ws.Range("C46").Value = "La Araucanía"
$ws.Range("D46").Value = 44669
$ws.Range("E46").Value = 9
$ws.Range("F46").Value = 100114002
$ws.Range("G46").Value = "Camote"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 25
$ws.Range("K46").Value = 18000
$ws.Range("L46").Value = 18000
$ws.Range("M46").Value = 18000
$ws.Range("N46").Value = "$/malla 20 kilos"
$ws.Range("O46").Value = "Perú"
$ws.Range("P46").Value = 900
$ws.Range("Q46").Value = 20
$ws.Range("R46").Value = "Hortaliza"
